$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, shifting rows 114:159 down to 115:160.
$ws.Rows("114:114").Insert()

# Populate the new row 114 with the new Puerro price record.
$ws.Cells.Item(114, 1).Value = 10
$ws.Cells.Item(114, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(114, 3).Value = "La Araucanía"
$ws.Cells.Item(114, 4).Value = 44523
$ws.Cells.Item(114, 5).Value = 9
$ws.Cells.Item(114, 6).Value = 100112005
$ws.Cells.Item(114, 7).Value = "Puerro"
$ws.Cells.Item(114, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 40
$ws.Cells.Item(114, 11).Value = 8000
$ws.Cells.Item(114, 12).Value = 8000
$ws.Cells.Item(114, 13).Value = 8000
$ws.Cells.Item(114, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(114, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(114, 16).Value = 667
$ws.Cells.Item(114, 17).Value = 12
$ws.Cells.Item(114, 18).Value = "Hortaliza"
